$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 4).Value = 3
$ws.Cells.Item(5, 4).Value = 89
$ws.Cells.Item(6, 4).Value = 3
$ws.Cells.Item(7, 4).Value = 13
$ws.Cells.Item(8, 4).Value = 5
$ws.Cells.Item(9, 4).Value = 5
$ws.Cells.Item(10, 4).Value = 5
$ws.Cells.Item(11, 4).Value = 5
$ws.Cells.Item(12, 4).Value = 13
$ws.Cells.Item(13, 4).Value = 3
$ws.Cells.Item(14, 4).Value = 5
$ws.Cells.Item(15, 4).Value = 5
$ws.Cells.Item(16, 4).Value = 5
$ws.Cells.Item(17, 4).Value = 3
$ws.Cells.Item(18, 4).Value = 3
$ws.Cells.Item(19, 4).Value = 3
$ws.Cells.Item(20, 4).Value = 89
$ws.Cells.Item(21, 4).Value = 13
$ws.Cells.Item(22, 4).Value = 13
$ws.Cells.Item(24, 4).Value = 8
$ws.Cells.Item(25, 4).Value = 5
$ws.Cells.Item(26, 4).Value = 21
$ws.Cells.Item(27, 4).Value = 8
$ws.Cells.Item(28, 4).Value = 8
$ws.Cells.Item(29, 4).Value = 8
$ws.Cells.Item(30, 4).Value = 8
$ws.Cells.Item(31, 4).Value = 8
$ws.Cells.Item(32, 4).Value = 8
$ws.Cells.Item(33, 4).Value = 8
$ws.Cells.Item(34, 4).Value = 8
$ws.Cells.Item(35, 4).Value = 8
$ws.Cells.Item(36, 4).Value = 8
$ws.Cells.Item(37, 4).Value = 8
$ws.Cells.Item(38, 4).Value = 8
$ws.Cells.Item(39, 4).Value = 8
$ws.Cells.Item(40, 4).Value = 8
$ws.Cells.Item(41, 4).Value = 8
$ws.Cells.Item(42, 4).Value = 8
$ws.Cells.Item(43, 4).Value = 8
$ws.Cells.Item(44, 4).Value = 8
$ws.Cells.Item(45, 4).Value = 8
$ws.Cells.Item(46, 4).Value = 8
$ws.Cells.Item(47, 4).Value = 8
$ws.Cells.Item(48, 4).Value = 8
$ws.Cells.Item(49, 4).Value = 8
$ws.Cells.Item(50, 4).Value = 8
$ws.Cells.Item(51, 4).Value = 8
$ws.Cells.Item(52, 4).Value = 8
$ws.Cells.Item(53, 4).Value = 8
$ws.Cells.Item(54, 4).Value = 8
$ws.Cells.Item(55, 4).Value = 8
$ws.Cells.Item(56, 4).Value = 8
$ws.Cells.Item(57, 4).Value = 8
$ws.Cells.Item(58, 4).Value = 8
$ws.Cells.Item(59, 4).Value = 8
$ws.Cells.Item(60, 4).Value = 8
$ws.Cells.Item(61, 4).Value = 8
$ws.Cells.Item(62, 4).Value = 3
$ws.Cells.Item(63, 4).Value = 3
$ws.Cells.Item(64, 4).Value = 13
$ws.Cells.Item(65, 4).Value = 5
$ws.Cells.Item(66, 4).Value = 5
$ws.Cells.Item(67, 4).Value = 5
$ws.Cells.Item(68, 4).Value = 21
$ws.Cells.Item(69, 4).Value = 21
$ws.Cells.Item(70, 4).Value = 3
$ws.Cells.Item(71, 4).Value = 5
$ws.Cells.Item(72, 4).Value = 8
$ws.Cells.Item(73, 4).Value = 5
$ws.Cells.Item(74, 4).Value = 5
$ws.Cells.Item(75, 4).Value = 5
$ws.Cells.Item(77, 4).Value = 8
$ws.Cells.Item(78, 4).Value = 89
$ws.Cells.Item(79, 4).Value = 3
$ws.Cells.Item(80, 4).Value = 5
$ws.Cells.Item(81, 4).Value = 3
$ws.Cells.Item(82, 4).Value = 34
$ws.Cells.Item(83, 4).Value = 21
$ws.Cells.Item(84, 4).Value = 8
$ws.Cells.Item(85, 4).Value = 3
$ws.Cells.Item(86, 4).Value = 5
$ws.Cells.Item(87, 4).Value = 89
$ws.Cells.Item(88, 4).Value = 21
$ws.Cells.Item(89, 4).Value = 3
$ws.Cells.Item(90, 4).Value = 5
$ws.Cells.Item(91, 4).Value = 55
$ws.Cells.Item(92, 4).Value = 8
$ws.Cells.Item(93, 4).Value = 8
$ws.Cells.Item(94, 4).Value = 21
$ws.Cells.Item(95, 4).Value = 5
$ws.Cells.Item(96, 4).Value = 89
$ws.Cells.Item(97, 4).Value = 3
$ws.Cells.Item(98, 4).Value = 3
$ws.Cells.Item(99, 4).Value = 5
$ws.Cells.Item(100, 4).Value = 5
$ws.Cells.Item(101, 4).Value = 5
$ws.Cells.Item(102, 4).Value = 3
$ws.Cells.Item(103, 4).Value = 21
$ws.Cells.Item(104, 4).Value = 89
$ws.Cells.Item(105, 4).Value = 89
$ws.Cells.Item(106, 4).Value = 21
$ws.Cells.Item(107, 4).Value = 3
$ws.Cells.Item(108, 4).Value = 3
$ws.Cells.Item(109, 4).Value = 5
$ws.Cells.Item(110, 4).Value = 8
$ws.Cells.Item(111, 4).Value = 8
$ws.Cells.Item(112, 4).Value = 3
$ws.Cells.Item(113, 4).Value = 5
$ws.Cells.Item(114, 4).Value = 8
$ws.Cells.Item(115, 4).Value = 5
$ws.Cells.Item(116, 4).Value = 21
$ws.Cells.Item(117, 4).Value = 21
$ws.Cells.Item(118, 4).Value = 21
$ws.Cells.Item(119, 4).Value = 8
$ws.Cells.Item(120, 4).Value = 3
$ws.Cells.Item(121, 4).Value = 3
$ws.Cells.Item(122, 4).Value = 3
$ws.Cells.Item(123, 4).Value = 5
$ws.Cells.Item(124, 4).Value = 8
$ws.Cells.Item(125, 4).Value = 8
$ws.Cells.Item(126, 4).Value = 8
$ws.Cells.Item(127, 4).Value = 89
$ws.Cells.Item(128, 4).Value = 3
$ws.Cells.Item(129, 4).Value = 3
$ws.Cells.Item(130, 4).Value = 34
$ws.Cells.Item(131, 4).Value = 8
$ws.Cells.Item(132, 4).Value = 3
$ws.Cells.Item(133, 4).Value = 3
$ws.Cells.Item(134, 4).Value = 5
$ws.Cells.Item(135, 4).Value = 3
$ws.Cells.Item(136, 4).Value = 21
$ws.Cells.Item(137, 4).Value = 89
$ws.Cells.Item(138, 4).Value = 3
$ws.Cells.Item(139, 4).Value = 8
$ws.Cells.Item(140, 4).Value = 13
$ws.Cells.Item(141, 4).Value = 3
$ws.Cells.Item(142, 4).Value = 21
$ws.Cells.Item(143, 4).Value = 21
$ws.Cells.Item(144, 4).Value = 21
$ws.Cells.Item(145, 4).Value = 21
$ws.Cells.Item(146, 4).Value = 21
$ws.Cells.Item(147, 4).Value = 21
$ws.Cells.Item(148, 4).Value = 21
$ws.Cells.Item(149, 4).Value = 21
$ws.Cells.Item(150, 4).Value = 21
$ws.Cells.Item(151, 4).Value = 21
$ws.Cells.Item(152, 4).Value = 21
$ws.Cells.Item(153, 4).Value = 21
$ws.Cells.Item(154, 4).Value = 21
$ws.Cells.Item(155, 4).Value = 21
$ws.Cells.Item(156, 4).Value = 21
$ws.Cells.Item(157, 4).Value = 21
$ws.Cells.Item(158, 4).Value = 21
$ws.Cells.Item(159, 4).Value = 21
$ws.Cells.Item(160, 4).Value = 21
$ws.Cells.Item(161, 4).Value = 21
$ws.Cells.Item(162, 4).Value = 21
$ws.Cells.Item(163, 4).Value = 21
$ws.Cells.Item(164, 4).Value = 21
$ws.Cells.Item(165, 4).Value = 21
$ws.Cells.Item(166, 4).Value = 21
$ws.Cells.Item(167, 4).Value = 8
$ws.Cells.Item(168, 4).Value = 3
$ws.Cells.Item(169, 4).Value = 3
$ws.Cells.Item(170, 4).Value = 3
$ws.Cells.Item(171, 4).Value = 21
$ws.Cells.Item(172, 4).Value = 5
$ws.Cells.Item(173, 4).Value = 5
$ws.Cells.Item(174, 4).Value = 89
$ws.Cells.Item(175, 4).Value = 5
$ws.Cells.Item(176, 4).Value = 34
$ws.Cells.Item(177, 4).Value = 3
$ws.Cells.Item(178, 4).Value = 3
$ws.Cells.Item(179, 4).Value = 55
$ws.Cells.Item(180, 4).Value = 55
$ws.Cells.Item(181, 4).Value = 55
$ws.Cells.Item(182, 4).Value = 55
$ws.Cells.Item(183, 4).Value = 55
$ws.Cells.Item(184, 4).Value = 55
$ws.Cells.Item(185, 4).Value = 55
$ws.Cells.Item(186, 4).Value = 55
$ws.Cells.Item(187, 4).Value = 55
$ws.Cells.Item(188, 4).Value = 55
$ws.Cells.Item(189, 4).Value = 55
$ws.Cells.Item(190, 4).Value = 55
$ws.Cells.Item(191, 4).Value = 55
$ws.Cells.Item(192, 4).Value = 55
$ws.Cells.Item(193, 4).Value = 55
$ws.Cells.Item(194, 4).Value = 3
$ws.Cells.Item(195, 4).Value = 8
$ws.Cells.Item(196, 4).Value = 8
$ws.Cells.Item(197, 4).Value = 55
$ws.Cells.Item(198, 4).Value = 34
$ws.Cells.Item(199, 4).Value = 3
$ws.Cells.Item(200, 4).Value = 3
$ws.Cells.Item(201, 4).Value = 3
$ws.Cells.Item(202, 4).Value = 13
$ws.Cells.Item(203, 4).Value = 3
$ws.Cells.Item(204, 4).Value = 5
$ws.Cells.Item(205, 4).Value = 89
$ws.Cells.Item(206, 4).Value = 3
$ws.Cells.Item(207, 4).Value = 3
$ws.Cells.Item(208, 4).Value = 3
$ws.Cells.Item(209, 4).Value = 3
$ws.Cells.Item(210, 4).Value = 3
$ws.Cells.Item(211, 4).Value = 3
$ws.Cells.Item(212, 4).Value = 5
$ws.Cells.Item(213, 4).Value = 3
$ws.Cells.Item(214, 4).Value = 5
$ws.Cells.Item(215, 4).Value = 5
$ws.Cells.Item(216, 4).Value = 5
$ws.Cells.Item(217, 4).Value = 3
$ws.Cells.Item(218, 4).Value = 5
$ws.Cells.Item(219, 4).Value = 5
$ws.Cells.Item(220, 4).Value = 5
$ws.Cells.Item(221, 4).Value = 5
$ws.Cells.Item(222, 4).Value = 5
$ws.Cells.Item(223, 4).Value = 5
$ws.Cells.Item(224, 4).Value = 5
$ws.Cells.Item(225, 4).Value = 5
$ws.Cells.Item(226, 4).Value = 5
$ws.Cells.Item(227, 4).Value = 5
$ws.Cells.Item(228, 4).Value = 5
$ws.Cells.Item(229, 4).Value = 21
$ws.Cells.Item(230, 4).Value = 34
$ws.Cells.Item(231, 4).Value = 34
$ws.Cells.Item(232, 4).Value = 8
$ws.Cells.Item(233, 4).Value = 13
$ws.Cells.Item(234, 4).Value = 5
$ws.Cells.Item(235, 4).Value = 3
$ws.Cells.Item(236, 4).Value = 3
$ws.Cells.Item(237, 4).Value = 8
$ws.Cells.Item(238, 4).Value = 21
$ws.Cells.Item(239, 4).Value = 55
$ws.Cells.Item(240, 4).Value = 13
$ws.Cells.Item(241, 4).Value = 3
$ws.Cells.Item(242, 4).Value = 3
$ws.Cells.Item(243, 4).Value = 3
$ws.Cells.Item(244, 4).Value = 3
$ws.Cells.Item(245, 4).Value = 3
$ws.Cells.Item(246, 4).Value = 8
$ws.Cells.Item(247, 4).Value = 8
$ws.Cells.Item(248, 4).Value = 89
$ws.Cells.Item(249, 4).Value = 89
$ws.Cells.Item(250, 4).Value = 89
$ws.Cells.Item(251, 4).Value = 89
$ws.Cells.Item(252, 4).Value = 89
$ws.Cells.Item(253, 4).Value = 89
$ws.Cells.Item(254, 4).Value = 55
$ws.Cells.Item(255, 4).Value = 3
$ws.Cells.Item(256, 4).Value = 13
$ws.Cells.Item(257, 4).Value = 3
$ws.Cells.Item(258, 4).Value = 3
$ws.Cells.Item(259, 4).Value = 13
$ws.Cells.Item(260, 4).Value = 3
$ws.Cells.Item(261, 4).Value = 13
$ws.Cells.Item(262, 4).Value = 5
$ws.Cells.Item(263, 4).Value = 8
$ws.Cells.Item(264, 4).Value = 21
$ws.Cells.Item(265, 4).Value = 3
$ws.Cells.Item(266, 4).Value = 8
$ws.Cells.Item(267, 4).Value = 8
$ws.Cells.Item(268, 4).Value = 8
$ws.Cells.Item(269, 4).Value = 55
$ws.Cells.Item(270, 4).Value = 3
$ws.Cells.Item(271, 4).Value = 5
$ws.Cells.Item(272, 4).Value = 8
$ws.Cells.Item(273, 4).Value = 5
$ws.Cells.Item(274, 4).Value = 8
$ws.Cells.Item(275, 4).Value = 5
$ws.Cells.Item(276, 4).Value = 3
$ws.Cells.Item(277, 4).Value = 8
$ws.Cells.Item(278, 4).Value = 5
$ws.Cells.Item(279, 4).Value = 21
$ws.Cells.Item(280, 4).Value = 5
$ws.Cells.Item(281, 4).Value = 3
$ws.Cells.Item(282, 4).Value = 8
$ws.Cells.Item(283, 4).Value = 13
$ws.Cells.Item(285, 4).Value = 3
$ws.Cells.Item(286, 4).Value = 3
$ws.Cells.Item(287, 4).Value = 5
$ws.Cells.Item(288, 4).Value = 3
$ws.Cells.Item(289, 4).Value = 13
$ws.Cells.Item(290, 4).Value = 13
$ws.Cells.Item(291, 4).Value = 13
$ws.Cells.Item(292, 4).Value = 13
$ws.Cells.Item(293, 4).Value = 13
$ws.Cells.Item(294, 4).Value = 13
$ws.Cells.Item(295, 4).Value = 13
$ws.Cells.Item(296, 4).Value = 13
$ws.Cells.Item(297, 4).Value = 13
$ws.Cells.Item(298, 4).Value = 13
$ws.Cells.Item(299, 4).Value = 13
$ws.Cells.Item(300, 4).Value = 13
$ws.Cells.Item(301, 4).Value = 13
$ws.Cells.Item(302, 4).Value = 13
$ws.Cells.Item(303, 4).Value = 13
$ws.Cells.Item(304, 4).Value = 13
$ws.Cells.Item(305, 4).Value = 13
$ws.Cells.Item(306, 4).Value = 13
$ws.Cells.Item(307, 4).Value = 13
$ws.Cells.Item(308, 4).Value = 13
$ws.Cells.Item(309, 4).Value = 13
$ws.Cells.Item(310, 4).Value = 13
$ws.Cells.Item(311, 4).Value = 13
$ws.Cells.Item(312, 4).Value = 13
$ws.Cells.Item(313, 4).Value = 13
$ws.Cells.Item(314, 4).Value = 13
$ws.Cells.Item(315, 4).Value = 13
$ws.Cells.Item(316, 4).Value = 13
$ws.Cells.Item(317, 4).Value = 13
$ws.Cells.Item(318, 4).Value = 13
$ws.Cells.Item(319, 4).Value = 13
$ws.Cells.Item(320, 4).Value = 13
$ws.Cells.Item(321, 4).Value = 13
$ws.Cells.Item(322, 4).Value = 13
$ws.Cells.Item(323, 4).Value = 13
$ws.Cells.Item(324, 4).Value = 13
$ws.Cells.Item(325, 4).Value = 13
$ws.Cells.Item(326, 4).Value = 13
$ws.Cells.Item(327, 4).Value = 13
$ws.Cells.Item(328, 4).Value = 13
$ws.Cells.Item(329, 4).Value = 13
$ws.Cells.Item(330, 4).Value = 13
$ws.Cells.Item(331, 4).Value = 13
$ws.Cells.Item(332, 4).Value = 13
$ws.Cells.Item(333, 4).Value = 8
$ws.Cells.Item(334, 4).Value = 89
$ws.Cells.Item(335, 4).Value = 3
$ws.Cells.Item(336, 4).Value = 3
$ws.Cells.Item(337, 4).Value = 13
$ws.Cells.Item(338, 4).Value = 3
$ws.Cells.Item(339, 4).Value = 8
$ws.Cells.Item(340, 4).Value = 89
$ws.Cells.Item(341, 4).Value = 5
$ws.Cells.Item(342, 4).Value = 34
$ws.Cells.Item(343, 4).Value = 13
$ws.Cells.Item(344, 4).Value = 5
$ws.Cells.Item(345, 4).Value = 3
$ws.Cells.Item(346, 4).Value = 8
$ws.Cells.Item(347, 4).Value = 13
$ws.Cells.Item(348, 4).Value = 13
$ws.Cells.Item(349, 4).Value = 13
$ws.Cells.Item(350, 4).Value = 8
$ws.Cells.Item(351, 4).Value = 21
$ws.Cells.Item(352, 4).Value = 89
$ws.Cells.Item(353, 4).Value = 89
$ws.Cells.Item(354, 4).Value = 21
$ws.Cells.Item(355, 4).Value = 8
$ws.Cells.Item(356, 4).Value = 89
$ws.Cells.Item(357, 4).Value = 89
$ws.Cells.Item(358, 4).Value = 13
$ws.Cells.Item(359, 4).Value = 5
$ws.Cells.Item(360, 4).Value = 8
$ws.Cells.Item(361, 4).Value = 8
$ws.Cells.Item(362, 4).Value = 89
$ws.Cells.Item(363, 4).Value = 5
$ws.Cells.Item(364, 4).Value = 3
$ws.Cells.Item(365, 4).Value = 3
$ws.Cells.Item(366, 4).Value = 3
$ws.Cells.Item(367, 4).Value = 5
$ws.Cells.Item(368, 4).Value = 3
$ws.Cells.Item(369, 4).Value = 8
$ws.Cells.Item(370, 4).Value = 8
$ws.Cells.Item(371, 4).Value = 8
$ws.Cells.Item(372, 4).Value = 3
$ws.Cells.Item(373, 4).Value = 3
$ws.Cells.Item(374, 4).Value = 5
$ws.Cells.Item(375, 4).Value = 5
$ws.Cells.Item(376, 4).Value = 5
$ws.Cells.Item(377, 4).Value = 34
$ws.Cells.Item(378, 4).Value = 3
$ws.Cells.Item(379, 4).Value = 5
$ws.Cells.Item(380, 4).Value = 5
$ws.Cells.Item(381, 4).Value = 5
$ws.Cells.Item(382, 4).Value = 5
$ws.Cells.Item(383, 4).Value = 5
$ws.Cells.Item(384, 4).Value = 3
$ws.Cells.Item(385, 4).Value = 3
$ws.Cells.Item(386, 4).Value = 8
$ws.Cells.Item(387, 4).Value = 3
$ws.Cells.Item(388, 4).Value = 21
$ws.Cells.Item(389, 4).Value = 3
$ws.Cells.Item(390, 4).Value = 3
$ws.Cells.Item(391, 4).Value = 3
$ws.Cells.Item(392, 4).Value = 13
$ws.Cells.Item(393, 4).Value = 21
$ws.Cells.Item(394, 4).Value = 5
$ws.Cells.Item(395, 4).Value = 13
$ws.Cells.Item(396, 4).Value = 3
$ws.Cells.Item(397, 4).Value = 3
$ws.Cells.Item(398, 4).Value = 3
$ws.Cells.Item(399, 4).Value = 3
$ws.Cells.Item(400, 4).Value = 3
$ws.Cells.Item(401, 4).Value = 3
$ws.Cells.Item(402, 4).Value = 3
$ws.Cells.Item(403, 4).Value = 3
$ws.Cells.Item(404, 4).Value = 3
$ws.Cells.Item(405, 4).Value = 3
$ws.Cells.Item(406, 4).Value = 3
$ws.Cells.Item(407, 4).Value = 3
$ws.Cells.Item(408, 4).Value = 3
$ws.Cells.Item(409, 4).Value = 3
$ws.Cells.Item(410, 4).Value = 3
$ws.Cells.Item(411, 4).Value = 3
$ws.Cells.Item(412, 4).Value = 3
$ws.Cells.Item(413, 4).Value = 3
$ws.Cells.Item(414, 4).Value = 3
$ws.Cells.Item(415, 4).Value = 3
$ws.Cells.Item(416, 4).Value = 3
$ws.Cells.Item(417, 4).Value = 3
$ws.Cells.Item(418, 4).Value = 3
$ws.Cells.Item(419, 4).Value = 3
$ws.Cells.Item(420, 4).Value = 3
$ws.Cells.Item(421, 4).Value = 3
$ws.Cells.Item(422, 4).Value = 3
$ws.Cells.Item(423, 4).Value = 3
$ws.Cells.Item(424, 4).Value = 3
$ws.Cells.Item(425, 4).Value = 3
$ws.Cells.Item(426, 4).Value = 3
$ws.Cells.Item(427, 4).Value = 3
$ws.Cells.Item(428, 4).Value = 3
$ws.Cells.Item(429, 4).Value = 3
$ws.Cells.Item(430, 4).Value = 3
$ws.Cells.Item(431, 4).Value = 3
$ws.Cells.Item(432, 4).Value = 3
$ws.Cells.Item(433, 4).Value = 3
$ws.Cells.Item(434, 4).Value = 3
$ws.Cells.Item(435, 4).Value = 3
$ws.Cells.Item(436, 4).Value = 3
$ws.Cells.Item(437, 4).Value = 3
$ws.Cells.Item(438, 4).Value = 3
$ws.Cells.Item(439, 4).Value = 3
$ws.Cells.Item(440, 4).Value = 3
$ws.Cells.Item(441, 4).Value = 3
$ws.Cells.Item(442, 4).Value = 3
$ws.Cells.Item(443, 4).Value = 3
$ws.Cells.Item(444, 4).Value = 3
$ws.Cells.Item(445, 4).Value = 3
$ws.Cells.Item(446, 4).Value = 3
$ws.Cells.Item(447, 4).Value = 3
$ws.Cells.Item(448, 4).Value = 3
$ws.Cells.Item(449, 4).Value = 3
$ws.Cells.Item(450, 4).Value = 3
$ws.Cells.Item(451, 4).Value = 3
$ws.Cells.Item(452, 4).Value = 3
$ws.Cells.Item(453, 4).Value = 3
$ws.Cells.Item(454, 4).Value = 3
$ws.Cells.Item(455, 4).Value = 3
$ws.Cells.Item(456, 4).Value = 3
$ws.Cells.Item(457, 4).Value = 3
$ws.Cells.Item(458, 4).Value = 3
$ws.Cells.Item(459, 4).Value = 3
$ws.Cells.Item(460, 4).Value = 3
$ws.Cells.Item(461, 4).Value = 3
$ws.Cells.Item(462, 4).Value = 3
$ws.Cells.Item(463, 4).Value = 3
$ws.Cells.Item(464, 4).Value = 3
$ws.Cells.Item(465, 4).Value = 3
$ws.Cells.Item(466, 4).Value = 3
$ws.Cells.Item(467, 4).Value = 3
$ws.Cells.Item(468, 4).Value = 3
$ws.Cells.Item(469, 4).Value = 3
$ws.Cells.Item(470, 4).Value = 3
$ws.Cells.Item(471, 4).Value = 5
$ws.Cells.Item(472, 4).Value = 5
$ws.Cells.Item(473, 4).Value = 5
$ws.Cells.Item(474, 4).Value = 3
$ws.Cells.Item(475, 4).Value = 5
$ws.Cells.Item(476, 4).Value = 13
$ws.Cells.Item(477, 4).Value = 3
$ws.Cells.Item(478, 4).Value = 34
$ws.Cells.Item(479, 4).Value = 8
$ws.Cells.Item(480, 4).Value = 3
$ws.Cells.Item(481, 4).Value = 5
$ws.Cells.Item(482, 4).Value = 3
$ws.Cells.Item(483, 4).Value = 5
$ws.Cells.Item(484, 4).Value = 3
$ws.Cells.Item(485, 4).Value = 3
$ws.Cells.Item(486, 4).Value = 3
$ws.Cells.Item(487, 4).Value = 8
$ws.Cells.Item(488, 4).Value = 21
$ws.Cells.Item(489, 4).Value = 13

Write-Output "Applied 483 changes"
